# "fixed duration contracts implemented"
#
# 1) Tweak a handful of existing numeric cells (floating point
#    re-computation noise + genuine value changes) on:
#      Production, RawMaterial, RawMaterialInventory, Contracts,
#      RawMaterialContract
# 2) Add two new sheets at the end of the workbook:
#      RawMaterialPrices, RawMaterialCosts_FD

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1a. Production
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Production")
$ws.Range("E3").Value  = 9.999999999999659
$ws.Range("E27").Value = 0.9999999999997726

# ---------------------------------------------------------------------------
# 1b. RawMaterial
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("RawMaterial")
$ws.Range("C2").Value  = 55.999999999999716
$ws.Range("C3").Value  = 0.0
$ws.Range("C4").Value  = 100.00000000000003
$ws.Range("C7").Value  = 10.000000000000002
$ws.Range("C12").Value = 8.5

# ---------------------------------------------------------------------------
# 1c. RawMaterialInventory
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("RawMaterialInventory")
$ws.Range("D3").Value  = 50.5
$ws.Range("D13").Value = 2.0
$ws.Range("D14").Value = 0.0
$ws.Range("D16").Value = 0.0
$ws.Range("D26").Value = 0.0

# ---------------------------------------------------------------------------
# 1d. Contracts
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Contracts")
$ws.Range("D2").Value  = 1.0
$ws.Range("D5").Value  = 0.0
$ws.Range("D8").Value  = 1.0
$ws.Range("D9").Value  = 0.0
$ws.Range("D24").Value = 1.0
$ws.Range("D25").Value = 0.0
$ws.Range("D38").Value = 1.0
$ws.Range("D41").Value = 0.0
$ws.Range("D42").Value = 1.0
$ws.Range("D45").Value = 0.0
$ws.Range("D46").Value = 0.0
$ws.Range("D48").Value = 1.0

# ---------------------------------------------------------------------------
# 1e. RawMaterialContract
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("RawMaterialContract")
$ws.Range("D2").Value  = 55.999999999999716
$ws.Range("D5").Value  = 0.0
$ws.Range("D9").Value  = 0.0
$ws.Range("D13").Value = 100.00000000000003
$ws.Range("D24").Value = 10.000000000000002
$ws.Range("D25").Value = 0.0
$ws.Range("D38").Value = 86.5
$ws.Range("D41").Value = 0.0
$ws.Range("D42").Value = 8.5
$ws.Range("D45").Value = 0.0

# ---------------------------------------------------------------------------
# 2a. New sheet: RawMaterialPrices (CALMONTH | RAW MATERIAL | METRIC)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsPrices = $wb.Worksheets.Add($null, $lastSheet)
$wsPrices.Name = "RawMaterialPrices"

$wsPrices.Cells.Item(1,1).Value = "CALMONTH"
$wsPrices.Cells.Item(1,2).Value = "RAW MATERIAL"
$wsPrices.Cells.Item(1,3).Value = "METRIC"

$pricesRows = @(
    @(202201, "R1", 0),
    @(202202, "R1", 10),
    @(202203, "R1", 10),
    @(202204, "R1", 10),
    @(202205, "R1", 10),
    @(202206, "R1", 10),
    @(202207, "R1", 10),
    @(202208, "R1", 10),
    @(202209, "R1", 10),
    @(202210, "R1", 2),
    @(202211, "R1", 2),
    @(202212, "R1", 56)
)

$r = 2
foreach ($row in $pricesRows) {
    $wsPrices.Cells.Item($r,1).Value = $row[0]
    $wsPrices.Cells.Item($r,2).Value = $row[1]
    $wsPrices.Cells.Item($r,3).Value = $row[2]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2b. New sheet: RawMaterialCosts_FD (RawMaterial | Period | Amount)
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsCosts = $wb.Worksheets.Add($null, $lastSheet)
$wsCosts.Name = "RawMaterialCosts_FD"

$wsCosts.Cells.Item(1,1).Value = "RawMaterial"
$wsCosts.Cells.Item(1,2).Value = "Period"
$wsCosts.Cells.Item(1,3).Value = "Amount"

$costsRows = @(
    @(202201, 0.0),
    @(202202, 0.0),
    @(202203, 500.0000000000001),
    @(202204, 252.5),
    @(202205, 262.49999999999994),
    @(202206, 0.0),
    @(202207, 49.999999999960444),
    @(202208, 275.0),
    @(202209, 382.49999999999994),
    @(202210, 0.0),
    @(202211, 0.0),
    @(202212, 0.0),
    @(202213, 0.0),
    @(202214, 0.0),
    @(202215, 0.0)
)

$r = 2
foreach ($row in $costsRows) {
    $wsCosts.Cells.Item($r,1).Value = "R1"
    $wsCosts.Cells.Item($r,2).Value = $row[0]
    $wsCosts.Cells.Item($r,3).Value = $row[1]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Restore the originally-active sheet/tab selection.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Production").Activate()

Write-Output "done"
